# Update the workbook's build/version strings from the January 30 build
# timestamp to the February 02 build timestamp, across every sheet/cell
# where the version string appears.

$wb = $excel.ActiveWorkbook

$oldVersion = "mines - January 30 (built on January 30 2026 16.19.47 EST)"
$newVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $val = $cell.Value2
        if ($val -ne $null -and $val -is [string] -and $val.Contains($oldVersion)) {
            $cell.Value2 = $val.Replace($oldVersion, $newVersion)
        }
    }
}
